# Applies the "Anonimyzed fedcore" update:
#  - renames the "fedcore" column header to "approach" on both sheets
#  - normalizes "-0" values to plain "0"
#  - gives the now-empty C1/D1 (and F1/G1) header cells -- which sit inside
#    the merged B1:D1 / E1:G1 banner cell -- a plain top/bottom (and,
#    for the rightmost one, also right) thin border instead of the full
#    boxed border they inherited from that banner cell's style
#  - drops the stray empty G5 cell on the computational_comparison sheet

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# --- header-cell borders -------------------------------------------------
# Build the "top+bottom only" border pattern from scratch exactly once
# (on sheet1!C1), derive the "top+bottom+right" pattern from it by
# flipping a single extra edge on a copy of it, and then stamp both
# finished patterns onto every other cell that needs them purely via
# copy / paste-of-formats. Re-deriving either pattern edge-by-edge a
# second time would leave an unused, but permanently registered, style
# behind in the workbook's style table -- reusing already-built styles
# keeps the style table exactly as small as the target.

$topBottomCell = $ws1.Range("C1")        # no left, no right, top+bottom thin
$topBottomCell.Style = "Normal"
$topBottomCell.Borders.Item(7).LineStyle  = -4142   # xlEdgeLeft   -> none
$topBottomCell.Borders.Item(8).LineStyle  = 1       # xlEdgeTop    -> continuous
$topBottomCell.Borders.Item(9).LineStyle  = 1       # xlEdgeBottom -> continuous
$topBottomCell.Borders.Item(10).LineStyle = -4142   # xlEdgeRight  -> none

$topBottomRightCell = $ws1.Range("D1")   # no left, top+bottom+right thin
$topBottomCell.Copy()
$topBottomRightCell.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$topBottomRightCell.Borders.Item(10).LineStyle = 1   # xlEdgeRight -> continuous

$topBottomCell.Copy()
$ws2.Range("C1").PasteSpecial(-4122)
$ws2.Range("F1").PasteSpecial(-4122)

$topBottomRightCell.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- rename the "fedcore" column headers to "approach" -----------------
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- normalize negative-zero "change" values to plain zero -------------
$ws1.Range("D4").Value  = 0
$ws1.Range("D5").Value  = 0
$ws1.Range("D12").Value = 0

# --- drop the stray empty G5 cell --------------------------------------
$ws2.Range("G5").ClearContents()
